{"js": "// Update the date line and the 100 multiplication-table answers.\n//\n// The document body is one paragraph of title/date text followed by a\n// 20-row x 5-column table (100 cells, one paragraph each) -- 101\n// paragraphs total, in document order. We overwrite each paragraph's\n// text in place (preserving its run formatting) with the corresponding\n// replacement below. Replacements are matched strictly by position\n// (not by searching for the old text) because a couple of the old\n// values repeat verbatim at different positions but map to different\n// new values.\nconst newValues = [\n  \"2023-04-13 Thursday\", \"93\u00d752=4836\", \"10\u00d793=930\", \"97\u00d713=1261\", \"90\u00d743=3870\", \"97\u00d719=1843\",\n  \"63\u00d776=4788\", \"66\u00d774=4884\", \"16\u00d784=1344\", \"69\u00d721=1449\", \"71\u00d758=4118\", \"92\u00d728=2576\",\n  \"54\u00d717=918\", \"16\u00d771=1136\", \"93\u00d757=5301\", \"28\u00d716=448\", \"59\u00d793=5487\", \"26\u00d740=1040\",\n  \"69\u00d748=3312\", \"68\u00d755=3740\", \"33\u00d795=3135\", \"68\u00d711=748\", \"68\u00d743=2924\", \"16\u00d713=208\",\n  \"74\u00d788=6512\", \"84\u00d746=3864\", \"48\u00d740=1920\", \"99\u00d795=9405\", \"72\u00d749=3528\", \"98\u00d747=4606\",\n  \"53\u00d725=1325\", \"20\u00d776=1520\", \"13\u00d747=611\", \"96\u00d766=6336\", \"54\u00d711=594\", \"96\u00d735=3360\",\n  \"73\u00d754=3942\", \"63\u00d736=2268\", \"19\u00d759=1121\", \"55\u00d781=4455\", \"36\u00d793=3348\", \"45\u00d743=1935\",\n  \"45\u00d724=1080\", \"28\u00d761=1708\", \"78\u00d747=3666\", \"76\u00d730=2280\", \"36\u00d718=648\", \"81\u00d760=4860\",\n  \"36\u00d769=2484\", \"48\u00d723=1104\", \"52\u00d710=520\", \"58\u00d728=1624\", \"45\u00d713=585\", \"86\u00d773=6278\",\n  \"24\u00d741=984\", \"67\u00d718=1206\", \"100\u00d758=5800\", \"62\u00d714=868\", \"90\u00d712=1080\", \"33\u00d772=2376\",\n  \"96\u00d757=5472\", \"15\u00d746=690\", \"86\u00d713=1118\", \"41\u00d722=902\", \"55\u00d736=1980\", \"36\u00d718=648\",\n  \"55\u00d740=2200\", \"60\u00d732=1920\", \"48\u00d799=4752\", \"19\u00d748=912\", \"34\u00d7100=3400\", \"39\u00d759=2301\",\n  \"67\u00d799=6633\", \"74\u00d794=6956\", \"91\u00d738=3458\", \"27\u00d772=1944\", \"69\u00d758=4002\", \"100\u00d792=9200\",\n  \"10\u00d723=230\", \"29\u00d732=928\", \"18\u00d790=1620\", \"35\u00d761=2135\", \"82\u00d731=2542\", \"49\u00d732=1568\",\n  \"34\u00d711=374\", \"35\u00d785=2975\", \"96\u00d799=9504\", \"87\u00d737=3219\", \"25\u00d739=975\", \"18\u00d764=1152\",\n  \"83\u00d747=3901\", \"26\u00d770=1820\", \"57\u00d773=4161\", \"79\u00d761=4819\", \"34\u00d718=612\", \"94\u00d718=1692\",\n  \"56\u00d748=2688\", \"40\u00d773=2920\", \"34\u00d7100=3400\", \"28\u00d775=2100\", \"15\u00d723=345\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].getRange().insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 100 multiplication-table answers.\n# $d.Paragraphs enumerates every paragraph mark in the document, including\n# the empty row-end paragraph Word inserts after each 5-cell table row.\n# We keep only the \"real\" (non-empty-after-trimming-control-chars)\n# paragraphs, in document order, and set each Range.Text to the matching\n# replacement value -- preserving each run's existing formatting since we\n# only change the text content of the existing range.\n$d = $word.ActiveDocument\n\n$newValues = @(\n    \"2023-04-13 Thursday\",\n    \"93\u00d752=4836\",\n    \"10\u00d793=930\",\n    \"97\u00d713=1261\",\n    \"90\u00d743=3870\",\n    \"97\u00d719=1843\",\n    \"63\u00d776=4788\",\n    \"66\u00d774=4884\",\n    \"16\u00d784=1344\",\n    \"69\u00d721=1449\",\n    \"71\u00d758=4118\",\n    \"92\u00d728=2576\",\n    \"54\u00d717=918\",\n    \"16\u00d771=1136\",\n    \"93\u00d757=5301\",\n    \"28\u00d716=448\",\n    \"59\u00d793=5487\",\n    \"26\u00d740=1040\",\n    \"69\u00d748=3312\",\n    \"68\u00d755=3740\",\n    \"33\u00d795=3135\",\n    \"68\u00d711=748\",\n    \"68\u00d743=2924\",\n    \"16\u00d713=208\",\n    \"74\u00d788=6512\",\n    \"84\u00d746=3864\",\n    \"48\u00d740=1920\",\n    \"99\u00d795=9405\",\n    \"72\u00d749=3528\",\n    \"98\u00d747=4606\",\n    \"53\u00d725=1325\",\n    \"20\u00d776=1520\",\n    \"13\u00d747=611\",\n    \"96\u00d766=6336\",\n    \"54\u00d711=594\",\n    \"96\u00d735=3360\",\n    \"73\u00d754=3942\",\n    \"63\u00d736=2268\",\n    \"19\u00d759=1121\",\n    \"55\u00d781=4455\",\n    \"36\u00d793=3348\",\n    \"45\u00d743=1935\",\n    \"45\u00d724=1080\",\n    \"28\u00d761=1708\",\n    \"78\u00d747=3666\",\n    \"76\u00d730=2280\",\n    \"36\u00d718=648\",\n    \"81\u00d760=4860\",\n    \"36\u00d769=2484\",\n    \"48\u00d723=1104\",\n    \"52\u00d710=520\",\n    \"58\u00d728=1624\",\n    \"45\u00d713=585\",\n    \"86\u00d773=6278\",\n    \"24\u00d741=984\",\n    \"67\u00d718=1206\",\n    \"100\u00d758=5800\",\n    \"62\u00d714=868\",\n    \"90\u00d712=1080\",\n    \"33\u00d772=2376\",\n    \"96\u00d757=5472\",\n    \"15\u00d746=690\",\n    \"86\u00d713=1118\",\n    \"41\u00d722=902\",\n    \"55\u00d736=1980\",\n    \"36\u00d718=648\",\n    \"55\u00d740=2200\",\n    \"60\u00d732=1920\",\n    \"48\u00d799=4752\",\n    \"19\u00d748=912\",\n    \"34\u00d7100=3400\",\n    \"39\u00d759=2301\",\n    \"67\u00d799=6633\",\n    \"74\u00d794=6956\",\n    \"91\u00d738=3458\",\n    \"27\u00d772=1944\",\n    \"69\u00d758=4002\",\n    \"100\u00d792=9200\",\n    \"10\u00d723=230\",\n    \"29\u00d732=928\",\n    \"18\u00d790=1620\",\n    \"35\u00d761=2135\",\n    \"82\u00d731=2542\",\n    \"49\u00d732=1568\",\n    \"34\u00d711=374\",\n    \"35\u00d785=2975\",\n    \"96\u00d799=9504\",\n    \"87\u00d737=3219\",\n    \"25\u00d739=975\",\n    \"18\u00d764=1152\",\n    \"83\u00d747=3901\",\n    \"26\u00d770=1820\",\n    \"57\u00d773=4161\",\n    \"79\u00d761=4819\",\n    \"34\u00d718=612\",\n    \"94\u00d718=1692\",\n    \"56\u00d748=2688\",\n    \"40\u00d773=2920\",\n    \"34\u00d7100=3400\",\n    \"28\u00d775=2100\",\n    \"15\u00d723=345\"\n)\n\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $t = $r.Text\n    $clean = $t -replace '[\\x00-\\x1F\\x7F]+$', ''\n    if ($clean.Length -gt 0) {\n        if ($i -ge $newValues.Count) {\n            throw \"More non-empty paragraphs ($($i+1)) than replacement values ($($newValues.Count))\"\n        }\n        $r.Text = $newValues[$i]\n        $i++\n    }\n}\n\nif ($i -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) non-empty paragraphs, processed $i\"\n}\n"}
